$wb = $excel.ActiveWorkbook

# Sheet "linear"
$wsLinear = $wb.Worksheets.Item("linear")
$wsLinear.Range("B2").Value = 0.0004954779980738566
$wsLinear.Range("B3").Value = -0.0009015378188562673
$wsLinear.Range("B4").Value = 0.0001299103179096192

# Sheet "non-linear"
$wsNonLinear = $wb.Worksheets.Item("non-linear")
$wsNonLinear.Range("B2").Value = -0.0003688182114076711
$wsNonLinear.Range("B3").Value = -0.002796682407822521
$wsNonLinear.Range("B4").Value = 0.0001869531197300377
$wsNonLinear.Range("B5").Value = 0.0004345906523447599
$wsNonLinear.Range("B6").Value = -0.0005255424046677973
$wsNonLinear.Range("B7").Value = 0.00008863470674000043
